$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# --- Rows 36/37 and 44/45 swap coin identity (Coin/Link/Price/Volume). ---
Set-TextValue 36 2 "Hedera"
Set-TextValue 36 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue 37 2 "WEMIXTOKEN"
Set-TextValue 37 3 "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue 44 2 "TheSandbox"
Set-TextValue 44 3 "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue 45 2 "EnergySwap"
Set-TextValue 45 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"

# --- Price (D) and Volume(1h) (E) updates. ---
Set-TextValue 2 4 "19.859.27"
Set-TextValue 2 5 "  -8.54%  "
Set-TextValue 3 4 "1.404.17"
Set-TextValue 3 5 "  -8.63%  "
Set-TextValue 4 4 "1.004"
Set-TextValue 4 5 "  +0.32%  "
Set-TextValue 5 4 "1.003"
Set-TextValue 5 5 "  +0.23%  "
Set-TextValue 6 4 "273.07"
Set-TextValue 6 5 "  -5.73%  "
Set-TextValue 7 4 "0.3673"
Set-TextValue 7 5 "  -6.47%  "
Set-TextValue 8 4 "0.3110"
Set-TextValue 8 5 "  -2.57%  "
Set-TextValue 9 4 "39.34"
Set-TextValue 9 5 "  -9.40%  "
Set-TextValue 10 4 "1.003"
Set-TextValue 10 5 "  -6.39%  "
Set-TextValue 11 4 "0.06475"
Set-TextValue 11 5 "  -9.89%  "
Set-TextValue 12 4 "1.004"
Set-TextValue 12 5 "  +0.30%  "
Set-TextValue 13 4 "5.420"
Set-TextValue 13 5 "  -5.76%  "
Set-TextValue 14 4 "17.39"
Set-TextValue 14 5 "  -5.64%  "
Set-TextValue 15 4 "6.147"
Set-TextValue 15 5 "  -7.24%  "
Set-TextValue 16 4 "1.404.24"
Set-TextValue 16 5 "  -8.48%  "
Set-TextValue 17 4 "0.00001010"
Set-TextValue 17 5 "  -7.73%  "
Set-TextValue 18 4 "0.05666"
Set-TextValue 18 5 "  -14.15%  "
Set-TextValue 19 4 "1.003"
Set-TextValue 19 5 "  +0.28%  "
Set-TextValue 20 4 "70.19"
Set-TextValue 20 5 "  -16.49%  "
Set-TextValue 21 4 "5.565"
Set-TextValue 21 5 "  -9.39%  "
Set-TextValue 22 4 "14.64"
Set-TextValue 22 5 "  -5.70%  "
Set-TextValue 23 5 "  +1.58%  "
Set-TextValue 24 4 "2.276"
Set-TextValue 24 5 "  -3.66%  "
Set-TextValue 25 4 "19.887.82"
Set-TextValue 25 5 "  -8.36%  "
Set-TextValue 26 4 "2.234"
Set-TextValue 26 5 "  -6.13%  "
Set-TextValue 27 4 "135.14"
Set-TextValue 27 5 "  -11.13%  "
Set-TextValue 28 4 "16.87"
Set-TextValue 28 5 "  -8.66%  "
Set-TextValue 29 4 "1.560.78"
Set-TextValue 29 5 "  -7.23%  "
Set-TextValue 30 4 "109.05"
Set-TextValue 30 5 "  -7.29%  "
Set-TextValue 31 4 "4.099"
Set-TextValue 31 5 "  -15.81%  "
Set-TextValue 32 4 "5.291"
Set-TextValue 32 5 "  -12.94%  "
Set-TextValue 33 4 "0.8093"
Set-TextValue 33 5 "  -16.03%  "
Set-TextValue 34 4 "0.07661"
Set-TextValue 34 5 "  -5.34%  "
Set-TextValue 35 4 "8.378"
Set-TextValue 35 5 "  -1.85%  "
Set-TextValue 36 4 "0.05851"
Set-TextValue 36 5 "  -1.97%  "
Set-TextValue 37 4 "1.440"
Set-TextValue 37 5 "  -4.11%  "
Set-TextValue 38 4 "4.804"
Set-TextValue 38 5 "  -7.44%  "
Set-TextValue 39 4 "1.002"
Set-TextValue 39 5 "  +0.23%  "
Set-TextValue 40 4 "0.02057"
Set-TextValue 40 5 "  -7.63%  "
Set-TextValue 41 4 "0.1899"
Set-TextValue 41 5 "  -7.08%  "
Set-TextValue 42 4 "10.32"
Set-TextValue 42 5 "  -8.61%  "
Set-TextValue 43 4 "1.083"
Set-TextValue 43 5 "  -8.24%  "
Set-TextValue 44 4 "0.5261"
Set-TextValue 44 5 "  -9.56%  "
Set-TextValue 45 4 "12.24"
Set-TextValue 45 5 "  -7.02%  "
Set-TextValue 46 4 "3.505"
Set-TextValue 46 5 "  -5.90%  "
Set-TextValue 47 4 "0.5101"
Set-TextValue 47 5 "  -8.66%  "
Set-TextValue 48 4 "111.59"
Set-TextValue 48 5 "  -3.83%  "
Set-TextValue 49 4 "1.758"
Set-TextValue 49 5 "  -7.13%  "
Set-TextValue 50 4 "1.034"
Set-TextValue 50 5 "  -11.29%  "
Set-TextValue 51 5 "  +0.20%  "
